$wb = $excel.ActiveWorkbook

# --- Sheet "result 1" ---
$ws1 = $wb.Worksheets.Item("result 1")
$ws1.Range("B2").Value = "NB13"
$ws1.Range("C2").Value = 11626
$ws1.Range("B3").Value = "PC20"
$ws1.Range("C3").Value = 7862

# --- Sheet "result 2" ---
$ws2 = $wb.Worksheets.Item("result 2")
$ws2.Range("B2").Value = 11626
$ws2.Range("C2").Value = 8675
$ws2.Range("D2").Value = 9728.33
$ws2.Range("B3").Value = 7862
$ws2.Range("C3").Value = 5854
$ws2.Range("D3").Value = 6858

# --- Sheet "result 3" ---
$ws3 = $wb.Worksheets.Item("result 3")
$ws3.Range("B2").Value = 1506
$ws3.Range("C2").Value = 1136
$ws3.Range("D2").Value = 1291
